$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Copy formatting (style) from the row above for the two date-formatted columns (A and G)
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("G5").Copy($ws.Range("G6"))

# Now set the actual values for the new row
$ws.Cells.Item($row, 1).Value = 42636.593182870369
$ws.Cells.Item($row, 2).Value = $false
$ws.Cells.Item($row, 3).Value = 9923.64
$ws.Cells.Item($row, 4).Value = 9951.5
$ws.Cells.Item($row, 5).Value = 313.07
$ws.Cells.Item($row, 6).Value = 314.81
$ws.Cells.Item($row, 7).Value = $true
$ws.Cells.Item($row, 8).Value = 0.56
$ws.Cells.Item($row, 9).Value = $false
